$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two empty homework-score cells for header row 2
$ws.Range("H2").Value = 6
$ws.Range("I2").Value = 7

# Fill in the two empty homework-score cells for row 10 (J10's SUM formula
# recalculates automatically: 25 -> 35)
$ws.Range("H10").Value = 5
$ws.Range("I10").Value = 5

# Update the view: zoom level and the active cell/selection in the
# frozen bottom-right pane
$ws.Application.ActiveWindow.Zoom = 115
$ws.Range("H10").Select()
